$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset entirely:
# original row 26 ("RM 232") and original row 28 ("SC 92").
# Deleting row 26 first shifts "SC 92" up to row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply the remaining cell-level value changes (re-imputed / re-masked values).
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = 17.64
$ws.Range("D5").Value = ""
$ws.Range("E6").Value = -5.7
$ws.Range("D8").Value = -13.9
$ws.Range("D10").Value = -14.7
$ws.Range("F10").Value = ""
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("E13").Value = -5.3
$ws.Range("D15").Value = -15.2
$ws.Range("F16").Value = 17.34
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = 17.78
$ws.Range("D18").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F20").Value = 17.73
$ws.Range("F24").Value = ""
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = ""
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = 17.39
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("F33").Value = 17.53
